$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new "Shortest Job First" output section right before the
#    "## Version History" paragraph.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*## Version History*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find '## Version History' paragraph"
}

$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)

$bigBlockXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p><w:r>
        <w:t>Shortest Job First Algorithm</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>:</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t>Real Time CPU Scheduling Simulator:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Menu:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>1. First Come First Serve CPU Scheduling Algorithm</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>2. Shortest Job First (Non-Preemptive) CPU Scheduling Algorithm</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>3. Shortest Job First (Preemptive) CPU Scheduling Algorithm</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>4. Round Robin CPU Scheduling Algorithm</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Enter option and press enter: 2</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>User selected option: 2</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Enter the number of processes: 3</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Shortest Job First CPU Scheduling Algorithm</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Enter process related information:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Process 1 id: 101</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 1 burst time in milliseconds: 10</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 2 id: 102</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 2 burst time in milliseconds: 6</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 3 id: 103</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 3 burst time in milliseconds: 4</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>User entered process related information is as follows:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Process 1 id: 101</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 1 burst time in milliseconds: 10</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 2 id: 102</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>Process 2 burst time in milliseconds: 6</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 3 id: 103</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process 3 burst time in milliseconds: 4</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Process Schedule, Waiting time &amp; turnaround time:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Process id: 103</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Waiting time in milliseconds: 0</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Turnaround time in milliseconds: 4</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Process id: 102</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Waiting time in milliseconds: 4</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Turnaround time in milliseconds: 10</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Process id: 101</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Waiting time in milliseconds: 10</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Turnaround time in milliseconds: 20</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Average </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>Waiting</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> time in milliseconds: 4.66667</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Average Turnaround time in milliseconds: 11.3333</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>-----------------------------------------------------------------</w:t>
      </w:r>
    </w:p><w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($bigBlockXml)

# ---------------------------------------------------------------------------
# 2. Add a <w:lastRenderedPageBreak/> before "* 1.0" (new page-break location
#    after the content that was just inserted).
# ---------------------------------------------------------------------------
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "* 1.0") {
        $target2 = $p
        break
    }
}

if ($target2 -eq $null) {
    throw "Could not find '* 1.0' paragraph"
}

$brkPoint = $d.Range($target2.Range.Start, $target2.Range.Start)

$brkXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:lastRenderedPageBreak/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$brkPoint.InsertXML($brkXml)

# ---------------------------------------------------------------------------
# 3. Remove the stale <w:lastRenderedPageBreak/> that used to precede
#    "This project is licensed under ..." (the page break moved earlier in
#    the document now that new content was inserted above).
# ---------------------------------------------------------------------------
$target3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This project is licensed*") {
        $target3 = $p
        break
    }
}

if ($target3 -eq $null) {
    throw "Could not find 'This project is licensed' paragraph"
}

$pStart = $target3.Range.Start
$pEndNoMark = $target3.Range.End - 1
$replaceRange = $d.Range($pStart, $pEndNoMark)

$licenseXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">This project is licensed under </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>the an</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> open source/free software License.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$replaceRange.InsertXML($licenseXml)

Write-Output "All edits applied"
